$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force Text format on the Price (D) and Volume (E) columns so that
# numeric-looking strings (e.g. "20.00", "0.9988", "30.177.51") are
# preserved exactly as text instead of being auto-converted to numbers.
$ws.Range("D2:E51").NumberFormat = "@"

$ws.Range('D2').Value = '30.177.51'
$ws.Range('E2').Value = '  -0.90%  '
$ws.Range('D3').Value = '1.858.41'
$ws.Range('E3').Value = '  -2.73%  '
$ws.Range('D4').Value = '0.9988'
$ws.Range('E4').Value = '  -0.03%  '
$ws.Range('D5').Value = '234.03'
$ws.Range('E5').Value = '  -2.42%  '
$ws.Range('D6').Value = '0.9995'
$ws.Range('D7').Value = '0.4695'
$ws.Range('E7').Value = '  -1.70%  '
$ws.Range('D8').Value = '0.2814'
$ws.Range('E8').Value = '  -0.82%  '
$ws.Range('D9').Value = '0.06553'
$ws.Range('E9').Value = '  -2.10%  '
$ws.Range('D10').Value = '20.00'
$ws.Range('E10').Value = '  +3.35%  '
$ws.Range('D11').Value = '0.07792'
$ws.Range('E11').Value = '  +0.16%  '
$ws.Range('D12').Value = '97.07'
$ws.Range('E12').Value = '  -5.45%  '
$ws.Range('D13').Value = '1.860.68'
$ws.Range('E13').Value = '  -3.35%  '
$ws.Range('D14').Value = '5.105'
$ws.Range('E14').Value = '  -1.94%  '
$ws.Range('D15').Value = '0.6651'
$ws.Range('E15').Value = '  -0.86%  '
$ws.Range('D16').Value = '282.84'
$ws.Range('E16').Value = '  -1.06%  '
$ws.Range('D17').Value = '30.181.58'
$ws.Range('E17').Value = '  -1.01%  '
$ws.Range('E18').Value = '  +0.10%  '
$ws.Range('D19').Value = '5.474'
$ws.Range('E19').Value = '  +1.61%  '
$ws.Range('D20').Value = '12.57'
$ws.Range('E20').Value = '  -0.59%  '
$ws.Range('D21').Value = '2.101.77'
$ws.Range('E21').Value = '  -2.69%  '
$ws.Range('B22').Value = 'ShibaInu'
$ws.Range('C22').Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range('D22').Value = '0.000007232'
$ws.Range('E22').Value = '  -3.26%  '
$ws.Range('B23').Value = 'BinanceUSD'
$ws.Range('C23').Value = 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'
$ws.Range('D23').Value = '0.9988'
$ws.Range('E23').Value = '  +0.01%  '
$ws.Range('D24').Value = '6.144'
$ws.Range('E24').Value = '  -2.50%  '
$ws.Range('D25').Value = '168.08'
$ws.Range('E25').Value = '  +0.43%  '
$ws.Range('D26').Value = '9.299'
$ws.Range('E26').Value = '  -0.72%  '
$ws.Range('D27').Value = '19.02'
$ws.Range('E27').Value = '  -1.63%  '
$ws.Range('D28').Value = '1.917'
$ws.Range('E28').Value = '  -8.04%  '
$ws.Range('D29').Value = '1.340'
$ws.Range('E29').Value = '  -2.98%  '
$ws.Range('D30').Value = '0.09593'
$ws.Range('E30').Value = '  -3.46%  '
$ws.Range('D31').Value = '4.424'
$ws.Range('E31').Value = '  -3.24%  '
$ws.Range('D32').Value = '1.469'
$ws.Range('E32').Value = '  -3.02%  '
$ws.Range('D33').Value = '4.103'
$ws.Range('E33').Value = '  -3.43%  '
$ws.Range('D34').Value = '0.04677'
$ws.Range('E34').Value = '  -1.21%  '
$ws.Range('D35').Value = '1.100'
$ws.Range('E35').Value = '  -0.80%  '
$ws.Range('D36').Value = '0.6982'
$ws.Range('E36').Value = '  -3.64%  '
$ws.Range('D37').Value = '0.9999'
$ws.Range('E37').Value = '  +0.08%  '
$ws.Range('D38').Value = '2.714'
$ws.Range('E38').Value = '  +0.25%  '
$ws.Range('D39').Value = '0.01850'
$ws.Range('E39').Value = '  -2.81%  '
$ws.Range('D40').Value = '6.317'
$ws.Range('E40').Value = '  -0.21%  '
$ws.Range('E41').Value = '  -4.24%  '
$ws.Range('D42').Value = '71.91'
$ws.Range('E42').Value = '  -3.09%  '
$ws.Range('D43').Value = '0.8569'
$ws.Range('E43').Value = '  -0.95%  '
$ws.Range('D44').Value = '1.942'
$ws.Range('E44').Value = '  -0.48%  '
$ws.Range('B45').Value = 'Quant'
$ws.Range('C45').Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range('D45').Value = '104.17'
$ws.Range('E45').Value = '  -1.43%  '
$ws.Range('B46').Value = 'TheSandbox'
$ws.Range('C46').Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range('D46').Value = '0.4163'
$ws.Range('E46').Value = '  -2.02%  '
$ws.Range('B47').Value = 'PaxDollar'
$ws.Range('C47').Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$ws.Range('D47').Value = '1.0000'
$ws.Range('E47').Value = '  +0.05%  '
$ws.Range('D48').Value = '1.023.51'
$ws.Range('E48').Value = '  +4.89%  '
$ws.Range('D49').Value = '7.211'
$ws.Range('E49').Value = '  -2.59%  '
$ws.Range('D50').Value = '8.918'
$ws.Range('E50').Value = '  +2.17%  '
$ws.Range('D51').Value = '33.70'
$ws.Range('E51').Value = '  -2.60%  '
